$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Rename the sheet/tab and its SubCategory label
$ws.Name = "Regulations"
$ws.Range("B2").Value = "Regulations"
$ws.Range("B3").Value = "Regulations"

# --- Row 2: Securities and Exchange Board of India (Registrars to an Issue and Share Transfer Agents) Regulations, 2025 ---
$ws.Range("C2").Value = 2025
$ws.Range("D2").Value = "December"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2025-12-16"
$ws.Range("F2").Value = "Securities and Exchange Board of India (Registrars to an Issue and Share Transfer Agents) Regulations, 2025"
$ws.Range("G2").Value = "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1766125665223.pdf"
$ws.Range("H2").Value = "1766125665223.pdf"
$ws.Range("I2").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/SEBI/Regulations/2025/December/1766125665223.pdf"
$ws.Range("J2").Value = "`"This SEBI regulation applies to individuals or entities seeking to act as Registrars to an Issue and Share Transfer Agents (RTAs). The regulation requires obtaining a certificate of registration from the Board. Key compliance areas include listing, net worth requirements, conduct, books of accounts maintenance, and reporting. RTAs must appoint a compliance officer for monitoring compliance and grievance redressal, and are required to adhere to Acts, rules, regulations, and a Code of Conduct.`""

$embed1 = @"
1. (1) These regulations may be called the Securities and Exchange Board of India (Registrars to an Issue and Share Transfer Agents) Regulations, 2025.
(2) They shall come into force on the date of their publication in the Official Gazette.
2. (1) In these regulations, unless the context otherwise requires, the terms defined herein shall bear the meaning assigned to them below and their cognate expressions and variations shall be construed accordingly,-
(A) if its shares are listed on any recognised stock exchange, shall be construed with reference to the definition of control in terms of regulations framed under clause (h) of sub-section (2) of section 11 of the Act;
(B) if its shares are not listed on any recognised stock exchange, shall be construed with reference to the definition of control as provided in sub-section (27) of Section 2 of the Companies Act, 2013 (18 of 2013);
(2) Words and expressions used and not defined in these regulations but defined in the Act, the Companies Act, 2013, the Securities Contracts (Regulation) Act, 1956, the Depositories Act, 1996 or any rules or regulations made thereunder shall have the same meanings respectively assigned to them in those Acts, rules or regulations made thereunder or any statutory modification or re-enactment thereto, as the case may be.
REGISTRATION OF REGISTRAR TO AN ISSUE AND SHARE TRANSFER AGENT
Registration as Registrar to an Issue and Share Transfer Agent.
3. No person shall act as a Registrar to an Issue and Share Transfer Agent unless it has obtained a certificate of registration from the Board under these regulations.
Application for grant of certificate of registration.
4. (1) An application by a Registrar to an Issue and Share Transfer Agent for grant of a certificate of registration shall be made to the Board in Form A of Schedule I and shall be accompanied by a non- refundable application fee, as specified in Para 1 of Schedule II. (2) Notwithstanding anything contained in sub-regulation (1), any application made by a registrar to an issue or a share transfer agent prior to coming into force of these regulations containing such particulars as near thereto as mentioned in Form A of Schedule I shall be treated as an application made in pursuance of sub-regulation (1) and dealt with accordingly.
Application to conform to the requirements.
5. Subject to the provisions of sub-regulation (2) of regulation 4, any application, which is not complete in all respects and does not conform to the instructions specified in the form, shall be rejected:
Provided that, before rejecting any such application, the applicant shall be given an opportunity to remove, within thirty days of the date of receipt of the relevant communication from the Board, such objections as may be indicated by the Board:
Provided further that, the Board may, on sufficient reason being shown, extend the time for removal of objections by such further time, not exceeding thirty days, as the Board may consider fit to enable the applicant to remove such objections.
To furnish further information, clarification etc.
6. (1) The Board may require the applicant to furnish such further information or clarifications regarding matters relevant to its activity as registrar to an issue and share transfer agent to consider the application for grant of a certificate of registration. (3) The applicant or its authorised representative shall, if so required by the Board, appear before the Board for personal representation.
7. For considering the grant of certificate of registration to the applicant, the Board shall take into account all matters which it deems are relevant to the activities relating to registrar to an issue and share transfer agent and in particular the following, namely, whether- (a) the applicant is a body corporate;
with the securities markets or any other regulatory or enforcement agency;
8. For the purposes of determining whether an applicant is a fit and proper person the Board may take into account the criteria specified in Schedule II of the Securities and Exchange Board of India (Intermediaries) Regulations, 2008.
9. The net worth requirement referred to in clause (g) of regulation 7 shall not be less than fifty lakh rupees.
Explanation.— For the purpose of this regulation, "net worth" means the value of the paid up capital, the free reserves and the securities premium as disclosed in the books of accounts of the applicant at the time of making the application under sub-regulation (1) of regulation 4:
Provided that a registrar to an issue and share transfer agent, who was granted a certificate of registration prior to the commencement of Securities and Exchange Board of India (Registrars to an Issue and Share Transfer Agents) Regulations, 2025, shall raise its net worth to the said minimum within a period of eighteen months from such commencement.
10. (1) The Board, on being satisfied that the applicant fulfils the requirements specified in regulation 7, shall send an intimation to the applicant and on receipt of the payment of registration fees as specified in paragraph (2) of Schedule II, grant a certificate in Form B of Schedule I. (2) The certificate of registration granted under sub-regulation (1) shall be valid unless it is suspended or cancelled by the Board.
(3) The registrar to an issue or share transfer agent who has already been granted certificate of registration by the Board, prior to the commencement of the Securities and Exchange Board of India (Registrars to an Issue and Share Transfer Agents) Regulations, 2025 shall be deemed to have been granted a certificate of registration, in terms of sub-regulation (1).
(5) The grant of certificate of registration shall be subject to payment of fee specified in regulation 14.
11. (1) The certificate of registration granted under regulation 10 shall, inter alia, be subject to the following conditions, namely:-
the registrar to an issue and share transfer agent shall abide by the provisions of the Act and these regulations;
change in control in such manner as may be specified by the Board;
regulation 9 at all times during the period of validity of the certificate;
(2) Nothing contained in clause (b) of sub-regulation (1) shall affect the obligation to obtain a fresh registration under section 12 of the Act in cases where it is applicable.
12. (1) Where a registrar to an issue and share transfer agent is providing its services to unlisted companies, then, —
Explanation.— The registrars to an issue and share transfer agents serving mutual funds, asset management companies, portfolio management services and alternative investment funds shall remain under the regulatory purview of the Board.
(2) A registrar to an issue and share transfer agent offering services to unlisted companies prior to the commencement of the Securities and Exchange Board of India (Registrars to an Issue and Share Transfer Agents) Regulations, 2025 shall segregate such activities into a separate business unit within a period of eighteen months from the date of notification of these regulations.
Procedure where registration is not granted.
13. (1) After considering an application made under regulation 4, if the Board is of the opinion that a certificate should not be granted to the applicant, it may reject the application after giving the applicant a reasonable opportunity of being heard. (2) The decision of the Board to reject the application shall be communicated to the applicant within thirty days of such decision, stating therein the grounds on which the application has been rejected.
(3) The applicant may, being aggrieved by the decision of the Board under sub regulation (2) apply within a period of thirty days from the date of receipt of such intimation, to the Board for reconsideration of its decision.
(4) On receipt of the application made under sub-regulation (3), the Board shall reconsider its decision
and communicate its findings thereon as soon as possible in writing to the applicant.
Payment of fees and the consequences of failure to pay fees.
14. (1) Every applicant eligible for grant of registration, shall pay the fees in such manner and within the period specified in Schedule II:
Provided that the Board may on sufficient cause being shown permit the registrar to an issue and share transfer agent to pay such fees at any time before the expiry of six months from the date on which such fees becomes due.
(2) Where a registrar to an issue and share transfer agent fails to pay the fees as provided in sub- regulation (1), the Board may suspend the certificate, whereupon the registrar to an issue and share transfer agent shall cease to carry on any of his activity as a registrar to an issue and share transfer agent.
15. Every registrar to an issue and share transfer agent shall at all times abide by the Code of Conduct as specified in Schedule III.
Registrar to an issue and share transfer agent not to act as such for an associate.
16. No registrar to an issue and share transfer agent shall act as such for any issue of securities by a body corporate of which it is an associate.
Explanation.─ For the purpose of this regulation, a registrar to an issue and share transfer agent or the body corporate as the case may be shall be deemed to be an associate of the other where:-
body corporate or of the registrar to an issue and share transfer agent, as the case may be.
Explanation.— The term 'relative' shall have the same meaning as is assigned to it under sub- section (77) of section 2 of the Companies Act, 2013.
To maintain proper books of accounts and records, etc.
17. (1) Every Registrar to an Issue and Share Transfer Agent shall keep and maintain the following books of accounts and documents, namely: -
Every registrar to an issue and share transfer agent which handles issues shall maintain the
a. all the applications received from investors in respect of an issue;
b. all applications of investors rejected and reasons therefor;
c. basis of allotment of securities to the investors as finalised in consultation with the stock
terms and conditions of purchase of securities;
list of names of allottees and non-allottees of the securities;
g. refund orders dispatched to investors in respect of application monies received from them
h. such other records as may be specified by the Board.
(3) of a body corporate on whose behalf it is carrying on the activities as transfer agent namely:-
Every registrar to an issue and share transfer agent shall maintain the following records in respect
list of holders of securities of such body corporate;
the names of transferor and transferee and the dates of transfer of securities;
c. such other records as may be specified by the Board.
(4) books of accounts, records and documents are maintained.
Every registrar to an issue and share transfer agent shall intimate the Board the place where the
(5) Without prejudice to sub- regulation (1), every Registrar to an Issue and Share Transfer Agent shall, after the close of each financial year as soon as possible but not later than six months from the close of the said period furnish to the Board if so required copies of the balance sheet, profit and loss account, statement of net worth requirement and such other documents as may be required by the Board under this regulation.
18. Subject to provisions of any other law, the registrar to an issue and share transfer agent shall preserve the books of accounts and other records and documents maintained under regulation 17 for a minimum period of eight years.
19. (1) Every Registrar to an Issue and Share Transfer Agent shall appoint a compliance officer who shall be responsible for monitoring the compliance of the Act, rules and regulations, notifications, guidelines, instructions etc. issued by the Board or the Central Government and for redressal of investors' grievances. (2) The compliance officer shall immediately and independently report to the Board any non- compliance observed by him.
20. All 
"@
$ws.Range("K2").Value = $embed1

# --- Row 3: Securities and Exchange Board of India (Listing Obligations and Disclosure Requirements) (Sixth Amendment) Regulations, 2025 ---
$ws.Range("C3").Value = 2025
$ws.Range("D3").Value = "December"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2025-12-16"
$ws.Range("F3").Value = "Securities and Exchange Board of India (Listing Obligations and Disclosure Requirements) (Sixth Amendment) Regulations, 2025"
$ws.Range("G3").Value = "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1765885385474.pdf"
$ws.Range("H3").Value = "1765885385474.pdf"
$ws.Range("I3").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/SEBI/Regulations/2025/December/1765885385474.pdf"

$summary2 = @"
• SEBI has amended the Securities and Exchange Board of India (Listing Obligations and Disclosure Requirements) Regulations, 2015 by changing the term "Share Transfer Agent" to "Registrar to an Issue and Share Transfer Agent" throughout the regulations.
• This change affects the terms used in various sections, schedules, and sub-paragraphs of the regulations, impacting entities that act as registrars or share transfer agents for listed companies.
• The practical effect is that these entities must now be referred to consistently as "Registrar to an Issue and Share Transfer Agent" in all regulatory requirements, enhancing clarity and standardization.
"@
$ws.Range("J3").Value = $summary2

$embed2 = @"
They shall come into force on the date of their publication in the Official Gazette.
In the Securities and Exchange Board of India (Listing Obligations and Disclosure Requirements) Regulations, 2015 –
in the heading, the words "Share Transfer Agent" shall be substituted with the words "Registrar to an Issue and Share Transfer Agent";
in sub-regulation (1), the words, "share transfer agent" shall be substituted with the words "registrar to an issue and share transfer agent";
in the proviso to sub-regulation (1), the words "Category II share transfer agent" shall be substituted with the words "Registrar to an issue and share transfer agent";
in sub-regulation (4), the words "share transfer agent" wherever they appear, shall be substituted with the words "registrar to an issue and share transfer agent".
in clause (1), the words "Registrar to an Issue and/or Share Transfer Agent" shall be substituted with the words "Registrar to an Issue and Share Transfer Agent";
in clause (2), the words "Share Transfer Agent" shall be substituted with the words "Registrar to an Issue and Share Transfer Agent".
In Schedule II, in Part D, in paragraph B, sub-paragraph (3), the words and symbols "Registrar & Share Transfer Agent" shall be substituted with the words "Registrar to an Issue and Share Transfer Agent".
in Part A, in paragraph A, sub-paragraph (8), the words "share transfer agent" shall be substituted with the words "registrar to an issue and share transfer agent";
in Part B, in paragraph A, sub-paragraph (28), the words "Registrar and Share Transfer Agent" shall be substituted with the words "Registrar to an Issue and Share Transfer Agent";
in Part C, in paragraph A, sub-paragraph (4), clause (e), the words "registrar to an issue and/or share transfer agent" shall be substituted with the words "registrar and share transfer agent".
In Schedule VI, in paragraph A, the words "share transfer agent" shall be substituted with the words "registrar to an issue and share transfer agent".
In Schedule VII, in paragraph C, sub-paragraph (2), clause (f), sub-clause (ii), the words "Share Transfer Agent" shall be substituted with the words "Registrar to an Issue and Share Transfer Agent".
1. The Securities and Exchange Board of India (Listing Obligations and Disclosure Requirements) Regulations, 2015 were published in the Gazette of India on 2nd September 2015 vide No. SEBI/LAD-NRO/GN/2015- 16/013.
2. The Securities and Exchange Board of India (Listing Obligations and Disclosure Requirements) Regulations, 2015, were subsequently amended on:
a) December 22, 2015 by the Securities and Exchange Board of India (Listing Obligations and Disclosure
Requirements) (Amendment) Regulations, 2015 vide notification no. SEBI/LAD-NRO/GN/2015-16/27.
b) May 25, 2016 by Securities and Exchange Board of India (Listing Obligations and Disclosure Requirements) (Amendment) Regulations, 2016 vide notification no. SEBI/LAD-NRO/GN/ 2016- 17/001.
d) January 4, 2017 by Securities and Exchange Board of India (Listing Obligations and Disclosure
e) February 15, 2017 by Securities and Exchange Board of India (Listing Obligations and Disclosure Requirements) (Amendment) Regulations, 2017 vide notification no. SEBI/LAD/NRO/GN/2016-17/029.
f) March 6, 2017 by the Securities and Exchange Board of India (Payment of Fees and Mode of Payment) (Amendment) Regulations, 2017 vide Notification No. LAD-NRO/GN/2016- 17/037 read with March 29, 2017 by Securities and Exchange Board of India (Payment of Fees and Mode of Payment) (Amendment) Regulations, 2017 vide notification no. SEBI/LAD/NRO/GN/2016-17/38.
m) March 29, 2019 by the Securities and Exchange Board of India (Listing Obligations and Disclosure Requirements) (Amendment) Regulations, 2019 vide notification no. SEBI/LAD-NRO/GN/2019/07.
January 10, 2020 by the Securities and Exchange Board of India (Listing Obligations and Disclosure Requirements) (Amendment) Regulations, 2020, vide notification no. SEBI/ LAD-NRO/GN/2020/02.
s) April 17, 2020 by the Securities and Exchange Board of India (Regulatory Sandbox) (Amendment)
Regulations, 2020 vide no. SEBI/LAD-NRO/GN/2020/10.
v) January 8, 2021 by the Securities and Exchange Board of India (Listing Obligations and Disclosure Requirements) (Amendment) Regulations, 2021, vide notification no. SEBI/ LAD-NRO/GN/2021/02.
x) August 3, 2021 by the Securities and Exchange Board of India (Regulatory Sandbox) (Amendment)
Regulations, 2021 vide notification no. No. SEBI/LAD-NRO/GN/2021/30.
cc) January 24, 2022 by the Securities and Exchange Board of India (Listing Obligations and Disclosure Requirements) (Amendment) Regulations, 2022, vide notification No. SEBI/LAD-NRO/GN/2022/66.
ee) April 11, 2022 by the Securities and Exchange Board of India (Listing Obligations and Disclosure
ff) April 25, 2022 by the Securities and Exchange Board of India (Listing Obligations and Disclosure
gg) July 25, 2022 by the Securities and Exchange Board of India (Listing Obligations and Disclosure
hh) November 14, 2022 by the Securities and Exchange Board of India (Listing Obligations and Disclosure
ii) December 05, 2022 by the Securities and Exchange Board of India (Listing Obligations and Disclosure
jj) January 17, 2023 by the Securities and Exchange Board of India (Listing Obligations and Disclosure
Requirements) (Amendment) Regulations, 2023 vide notification No. SEBI/LAD-NRO/GN/- 2023/117.
kk) February 7, 2023 by the Securities and Exchange Board of India (Payment of Fees and Mode of
Payment) (Amendment) Regulations, 2023 vide notification No. SEBI/LAD-NRO/GN/2023/121.
July 4, 2023 by the Securities and Exchange Board of India (Alternative Dispute Resolution
Mechanism) (Amendment) Regulations, 2023 vide notification No. SEBI/LAD–NRO/GN/2023/137.
ss) May 17, 2024 by the Securities and Exchange Board of India (Listing Obligations and Disclosure Requirements) (Amendment) Regulations, 2024 vide notification No. SEBI/LAD-NRO/GN/2024/177.
vv) March 28, 2025 by the Securities and Exchange Board of India (Listing Obligations and Disclosure
Requirements) (Amendment) Regulations, 2025, vide notification No. SEBI/LAD-NRO/GN/2025/239.
Uploaded by Dte. of Printing at Government of India Press, Ring Road, Mayapuri, New Delhi-110064 and Published by the Controller of Publications, Delhi-110054.
"@
$ws.Range("K3").Value = $embed2

# Row heights as per diff
$ws.Range("A2:K2").RowHeight = 794.25
$ws.Range("A3:K3").RowHeight = 512.25
